$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2026-01-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-17 Saturday", 2)
$null = $d.Content.Find.Execute("94-62=32", $true, $false, $false, $false, $false, $true, 1, $false, "0+17=17", 2)
$null = $d.Content.Find.Execute("75+7=82", $true, $false, $false, $false, $false, $true, 1, $false, "62-48=14", 2)
$null = $d.Content.Find.Execute("23+6=29", $true, $false, $false, $false, $false, $true, 1, $false, "81-34=47", 2)
$null = $d.Content.Find.Execute("31+38=69", $true, $false, $false, $false, $false, $true, 1, $false, "1+62=63", 2)
$null = $d.Content.Find.Execute("97-65=32", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=95", 2)
$null = $d.Content.Find.Execute("33+63=96", $true, $false, $false, $false, $false, $true, 1, $false, "44-10=34", 2)
$null = $d.Content.Find.Execute("63+29=92", $true, $false, $false, $false, $false, $true, 1, $false, "66+19=85", 2)
$null = $d.Content.Find.Execute("18+23=41", $true, $false, $false, $false, $false, $true, 1, $false, "95-69=26", 2)
$null = $d.Content.Find.Execute("57+7=64", $true, $false, $false, $false, $false, $true, 1, $false, "86-75=11", 2)
$null = $d.Content.Find.Execute("91-61=30", $true, $false, $false, $false, $false, $true, 1, $false, "32+58=90", 2)
$null = $d.Content.Find.Execute("58-39=19", $true, $false, $false, $false, $false, $true, 1, $false, "53+6=59", 2)
$null = $d.Content.Find.Execute("50+22=72", $true, $false, $false, $false, $false, $true, 1, $false, "31-12=19", 2)
$null = $d.Content.Find.Execute("34+16=50", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=9", 2)
$null = $d.Content.Find.Execute("98-75=23", $true, $false, $false, $false, $false, $true, 1, $false, "1+69=70", 2)
$null = $d.Content.Find.Execute("0+71=71", $true, $false, $false, $false, $false, $true, 1, $false, "94-60=34", 2)
$null = $d.Content.Find.Execute("69+13=82", $true, $false, $false, $false, $false, $true, 1, $false, "39+14=53", 2)
$null = $d.Content.Find.Execute("33+58=91", $true, $false, $false, $false, $false, $true, 1, $false, "16+46=62", 2)
$null = $d.Content.Find.Execute("45+9=54", $true, $false, $false, $false, $false, $true, 1, $false, "86-55=31", 2)
$null = $d.Content.Find.Execute("9+15=24", $true, $false, $false, $false, $false, $true, 1, $false, "5+58=63", 2)
$null = $d.Content.Find.Execute("23+27=50", $true, $false, $false, $false, $false, $true, 1, $false, "10+63=73", 2)
$null = $d.Content.Find.Execute("50+20=70", $true, $false, $false, $false, $false, $true, 1, $false, "94-47=47", 2)
$null = $d.Content.Find.Execute("31+50=81", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=18", 2)
$null = $d.Content.Find.Execute("19+65=84", $true, $false, $false, $false, $false, $true, 1, $false, "23+41=64", 2)
$null = $d.Content.Find.Execute("16+39=55", $true, $false, $false, $false, $false, $true, 1, $false, "37+46=83", 2)
$null = $d.Content.Find.Execute("83-57=26", $true, $false, $false, $false, $false, $true, 1, $false, "55+36=91", 2)
$null = $d.Content.Find.Execute("99-72=27", $true, $false, $false, $false, $false, $true, 1, $false, "69-28=41", 2)
$null = $d.Content.Find.Execute("0+51=51", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=56", 2)
$null = $d.Content.Find.Execute("29+54=83", $true, $false, $false, $false, $false, $true, 1, $false, "28+60=88", 2)
$null = $d.Content.Find.Execute("26+25=51", $true, $false, $false, $false, $false, $true, 1, $false, "74-63=11", 2)
$null = $d.Content.Find.Execute("8+84=92", $true, $false, $false, $false, $false, $true, 1, $false, "17+64=81", 2)
$null = $d.Content.Find.Execute("31+62=93", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=48", 2)
$null = $d.Content.Find.Execute("86-28=58", $true, $false, $false, $false, $false, $true, 1, $false, "88-57=31", 2)
$null = $d.Content.Find.Execute("48-33=15", $true, $false, $false, $false, $false, $true, 1, $false, "32+20=52", 2)
$null = $d.Content.Find.Execute("29+47=76", $true, $false, $false, $false, $false, $true, 1, $false, "5+91=96", 2)
$null = $d.Content.Find.Execute("72-66=6", $true, $false, $false, $false, $false, $true, 1, $false, "6+19=25", 2)
$null = $d.Content.Find.Execute("32+13=45", $true, $false, $false, $false, $false, $true, 1, $false, "93-12=81", 2)
$null = $d.Content.Find.Execute("85-6=79", $true, $false, $false, $false, $false, $true, 1, $false, "45+18=63", 2)
$null = $d.Content.Find.Execute("54+0=54", $true, $false, $false, $false, $false, $true, 1, $false, "89-56=33", 2)
$null = $d.Content.Find.Execute("30-3=27", $true, $false, $false, $false, $false, $true, 1, $false, "62-21=41", 2)
$null = $d.Content.Find.Execute("80-48=32", $true, $false, $false, $false, $false, $true, 1, $false, "83-31=52", 2)
$null = $d.Content.Find.Execute("13+78=91", $true, $false, $false, $false, $false, $true, 1, $false, "62+20=82", 2)
$null = $d.Content.Find.Execute("46+3=49", $true, $false, $false, $false, $false, $true, 1, $false, "87-13=74", 2)
$null = $d.Content.Find.Execute("19-3=16", $true, $false, $false, $false, $false, $true, 1, $false, "86-85=1", 2)
$null = $d.Content.Find.Execute("99-82=17", $true, $false, $false, $false, $false, $true, 1, $false, "86-77=9", 2)
$null = $d.Content.Find.Execute("8+59=67", $true, $false, $false, $false, $false, $true, 1, $false, "51+14=65", 2)
$null = $d.Content.Find.Execute("64-38=26", $true, $false, $false, $false, $false, $true, 1, $false, "56+6=62", 2)
$null = $d.Content.Find.Execute("16-9=7", $true, $false, $false, $false, $false, $true, 1, $false, "94+3=97", 2)
$null = $d.Content.Find.Execute("19+16=35", $true, $false, $false, $false, $false, $true, 1, $false, "22+18=40", 2)
$null = $d.Content.Find.Execute("18+34=52", $true, $false, $false, $false, $false, $true, 1, $false, "5+71=76", 2)
$null = $d.Content.Find.Execute("25-19=6", $true, $false, $false, $false, $false, $true, 1, $false, "32+48=80", 2)
$null = $d.Content.Find.Execute("91-17=74", $true, $false, $false, $false, $false, $true, 1, $false, "82-55=27", 2)
$null = $d.Content.Find.Execute("27+10=37", $true, $false, $false, $false, $false, $true, 1, $false, "83-8=75", 2)
$null = $d.Content.Find.Execute("87-37=50", $true, $false, $false, $false, $false, $true, 1, $false, "86-25=61", 2)
$null = $d.Content.Find.Execute("50+16=66", $true, $false, $false, $false, $false, $true, 1, $false, "86-13=73", 2)
$null = $d.Content.Find.Execute("62-8=54", $true, $false, $false, $false, $false, $true, 1, $false, "36+29=65", 2)
$null = $d.Content.Find.Execute("25-2=23", $true, $false, $false, $false, $false, $true, 1, $false, "0+5=5", 2)
$null = $d.Content.Find.Execute("89-50=39", $true, $false, $false, $false, $false, $true, 1, $false, "5+37=42", 2)
$null = $d.Content.Find.Execute("67-21=46", $true, $false, $false, $false, $false, $true, 1, $false, "49-20=29", 2)
$null = $d.Content.Find.Execute("86-60=26", $true, $false, $false, $false, $false, $true, 1, $false, "80+1=81", 2)
$null = $d.Content.Find.Execute("47+6=53", $true, $false, $false, $false, $false, $true, 1, $false, "66+18=84", 2)
$null = $d.Content.Find.Execute("22-12=10", $true, $false, $false, $false, $false, $true, 1, $false, "48+12=60", 2)
$null = $d.Content.Find.Execute("65-52=13", $true, $false, $false, $false, $false, $true, 1, $false, "38+45=83", 2)
$null = $d.Content.Find.Execute("26-20=6", $true, $false, $false, $false, $false, $true, 1, $false, "22+20=42", 2)
$null = $d.Content.Find.Execute("56-12=44", $true, $false, $false, $false, $false, $true, 1, $false, "74-57=17", 2)
$null = $d.Content.Find.Execute("33-19=14", $true, $false, $false, $false, $false, $true, 1, $false, "76-16=60", 2)
$null = $d.Content.Find.Execute("24+23=47", $true, $false, $false, $false, $false, $true, 1, $false, "98-32=66", 2)
$null = $d.Content.Find.Execute("34+28=62", $true, $false, $false, $false, $false, $true, 1, $false, "69-67=2", 2)
$null = $d.Content.Find.Execute("83-76=7", $true, $false, $false, $false, $false, $true, 1, $false, "21+69=90", 2)
$null = $d.Content.Find.Execute("72+3=75", $true, $false, $false, $false, $false, $true, 1, $false, "63-52=11", 2)
$null = $d.Content.Find.Execute("68-39=29", $true, $false, $false, $false, $false, $true, 1, $false, "68-24=44", 2)
$null = $d.Content.Find.Execute("98-93=5", $true, $false, $false, $false, $false, $true, 1, $false, "1+98=99", 2)
$null = $d.Content.Find.Execute("65-12=53", $true, $false, $false, $false, $false, $true, 1, $false, "15+57=72", 2)
$null = $d.Content.Find.Execute("14+20=34", $true, $false, $false, $false, $false, $true, 1, $false, "0+2=2", 2)
$null = $d.Content.Find.Execute("4+29=33", $true, $false, $false, $false, $false, $true, 1, $false, "55-9=46", 2)
$null = $d.Content.Find.Execute("27-7=20", $true, $false, $false, $false, $false, $true, 1, $false, "30+51=81", 2)
$null = $d.Content.Find.Execute("54+11=65", $true, $false, $false, $false, $false, $true, 1, $false, "53+11=64", 2)
$null = $d.Content.Find.Execute("0+20=20", $true, $false, $false, $false, $false, $true, 1, $false, "52+24=76", 2)
$null = $d.Content.Find.Execute("40-18=22", $true, $false, $false, $false, $false, $true, 1, $false, "2+10=12", 2)
$null = $d.Content.Find.Execute("84-23=61", $true, $false, $false, $false, $false, $true, 1, $false, "12+6=18", 2)
$null = $d.Content.Find.Execute("55+11=66", $true, $false, $false, $false, $false, $true, 1, $false, "79-38=41", 2)
$null = $d.Content.Find.Execute("36-28=8", $true, $false, $false, $false, $false, $true, 1, $false, "47-47=0", 2)
$null = $d.Content.Find.Execute("92-50=42", $true, $false, $false, $false, $false, $true, 1, $false, "61-19=42", 2)
$null = $d.Content.Find.Execute("3+53=56", $true, $false, $false, $false, $false, $true, 1, $false, "20+70=90", 2)
$null = $d.Content.Find.Execute("34+24=58", $true, $false, $false, $false, $false, $true, 1, $false, "95-6=89", 2)
$null = $d.Content.Find.Execute("67-60=7", $true, $false, $false, $false, $false, $true, 1, $false, "31+44=75", 2)
$null = $d.Content.Find.Execute("70-45=25", $true, $false, $false, $false, $false, $true, 1, $false, "15+49=64", 2)
$null = $d.Content.Find.Execute("42-31=11", $true, $false, $false, $false, $false, $true, 1, $false, "3+57=60", 2)
$null = $d.Content.Find.Execute("25+28=53", $true, $false, $false, $false, $false, $true, 1, $false, "87-15=72", 2)
$null = $d.Content.Find.Execute("6+39=45", $true, $false, $false, $false, $false, $true, 1, $false, "86-6=80", 2)
$null = $d.Content.Find.Execute("25+51=76", $true, $false, $false, $false, $false, $true, 1, $false, "52-33=19", 2)
$null = $d.Content.Find.Execute("0+59=59", $true, $false, $false, $false, $false, $true, 1, $false, "61+29=90", 2)
$null = $d.Content.Find.Execute("6+20=26", $true, $false, $false, $false, $false, $true, 1, $false, "34-20=14", 2)
$null = $d.Content.Find.Execute("92-91=1", $true, $false, $false, $false, $false, $true, 1, $false, "91-23=68", 2)
$null = $d.Content.Find.Execute("42-28=14", $true, $false, $false, $false, $false, $true, 1, $false, "58+12=70", 2)
$null = $d.Content.Find.Execute("10+32=42", $true, $false, $false, $false, $false, $true, 1, $false, "73-1=72", 2)
$null = $d.Content.Find.Execute("64-27=37", $true, $false, $false, $false, $false, $true, 1, $false, "77-60=17", 2)
$null = $d.Content.Find.Execute("86+1=87", $true, $false, $false, $false, $false, $true, 1, $false, "73-15=58", 2)
$null = $d.Content.Find.Execute("4+76=80", $true, $false, $false, $false, $false, $true, 1, $false, "56-27=29", 2)
$null = $d.Content.Find.Execute("11+62=73", $true, $false, $false, $false, $false, $true, 1, $false, "23-16=7", 2)
$null = $d.Content.Find.Execute("40+7=47", $true, $false, $false, $false, $false, $true, 1, $false, "74+17=91", 2)
